# Apply update to moldova_super-liga_2023-2024 sheet:
#  1) Swap the content (columns F:V) of rows 40 and 41.
#  2) Append a new row 49 with a new match result, extending the
#     sheet dimension from A1:V48 to A1:V49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap rows 40 and 41 (columns F..V hold the match-specific data that was
#    reordered; A..E - Indice/pais/torneio/temporada/data_partida - stay put).
# ---------------------------------------------------------------------------

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$row40 = @()
$row41 = @()
foreach ($c in $cols) {
    $row40 += , ($ws.Range($c + "40").Value2)
    $row41 += , ($ws.Range($c + "41").Value2)
}

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "40").Value2 = $row41[$i]
    $ws.Range($cols[$i] + "41").Value2 = $row40[$i]
}

# ---------------------------------------------------------------------------
# 2) Append new row 49 (index 48) with the Floresti vs Sheriff Tiraspol match
#    played on 12/11/2023. Copy formatting from row 48 first so the new row
#    keeps the same per-column styles (bold/border on A, date format on E).
# ---------------------------------------------------------------------------

$ws.Range("A48:V48").Copy()
$ws.Range("A49:V49").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A49").Value2 = 48
$ws.Range("B49").Value2 = "moldova"
$ws.Range("C49").Value2 = "super-liga"
$ws.Range("D49").Value2 = "2023-2024"
$ws.Range("E49").Value2 = 45242.5
$ws.Range("F49").Value2 = "Floresti"
$ws.Range("G49").Value2 = 0
$ws.Range("H49").Value2 = "Sheriff Tiraspol"
$ws.Range("I49").Value2 = 4
$ws.Range("J49").Value2 = 28.55
$ws.Range("K49").Value2 = "12/11/2023 11:05"
$ws.Range("L49").Value2 = 28.55
$ws.Range("M49").Value2 = "12/11/2023 11:05"
$ws.Range("N49").Value2 = 16.03
$ws.Range("O49").Value2 = "12/11/2023 11:05"
$ws.Range("P49").Value2 = 16.03
$ws.Range("Q49").Value2 = "12/11/2023 11:05"
$ws.Range("R49").Value2 = 1.02
$ws.Range("S49").Value2 = "12/11/2023 10:58"
$ws.Range("T49").Value2 = 1.02
$ws.Range("U49").Value2 = "12/11/2023 10:58"
$ws.Range("V49").Value2 = "https://www.betexplorer.com/football/moldova/super-liga/floresti-sheriff-tiraspol/d4T2pQuF/"
